# Updates "想去人数" (F) and "最低票价" (G) numbers on sheets 展览 and 全部类型,
# and F2 on sheet 演出, per the upstream data refresh (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 13903
$ws1.Range("G5").Value  = 50
$ws1.Range("F6").Value  = 186
$ws1.Range("G6").Value  = 56
$ws1.Range("F7").Value  = 284
$ws1.Range("G7").Value  = 50
$ws1.Range("F8").Value  = 499
$ws1.Range("F9").Value  = 13
$ws1.Range("F14").Value = 463
$ws1.Range("F15").Value = 5909
$ws1.Range("F18").Value = 986
$ws1.Range("F19").Value = 128
$ws1.Range("F21").Value = 162
$ws1.Range("F22").Value = 286

# --- Sheet: 演出 ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 15

# --- Sheet: 全部类型 ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 13903
$ws4.Range("G5").Value  = 50
$ws4.Range("F6").Value  = 186
$ws4.Range("G6").Value  = 56
$ws4.Range("F7").Value  = 284
$ws4.Range("G7").Value  = 50
$ws4.Range("F8").Value  = 499
$ws4.Range("F9").Value  = 13
$ws4.Range("F14").Value = 463
$ws4.Range("F15").Value = 5909
$ws4.Range("F18").Value = 986
$ws4.Range("F19").Value = 128
$ws4.Range("F21").Value = 162
$ws4.Range("F22").Value = 286
$ws4.Range("F23").Value = 15
